$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text
# (matches the original sheet, where every cell - including the
# numeric-looking ones like totalRuns/totalBalls/sr - is typed as a
# string, t="str"). Setting NumberFormat to "@" before assigning the
# value stops Excel from re-interpreting "0", "1", "0.00" etc. as
# numbers; resetting the Style back to "Normal" afterwards clears the
# temporary text-format styling so the new row does not pick up an
# extra cell style compared with the existing rows.
function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# New row 3: another ball-by-ball / match record for Mitchell Marsh,
# duplicating the shape of the existing row 2 with the new match's data.
Set-TextCell 3 1  " Dubai (DSC)"
Set-TextCell 3 2  " September 21 2020"
Set-TextCell 3 3  "RCB won by 10 runs"
Set-TextCell 3 4  "Sunrisers Hyderabad"
Set-TextCell 3 5  "Royal Challengers Bangalore"
Set-TextCell 3 6  "Mitchell Marsh "
Set-TextCell 3 7  "0"
Set-TextCell 3 8  "1"
Set-TextCell 3 9  "0"
Set-TextCell 3 10 "0"
Set-TextCell 3 11 "0.00"
